$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Pass 1: seed the shared-strings table so new unique strings are created in
# the exact order the workbook ends up with (first-use order matters for the
# shared string table index). The values written here are overwritten with
# their real final values further down, so only the *order* matters.
# ---------------------------------------------------------------------------
$ws.Range("C13").Value = "Arsenal"
$ws.Range("D4").Value  = "Cardiff City"
$ws.Range("C2").Value  = "Manchester United"
$ws.Range("D2").Value  = "Burnley"
$ws.Range("C3").Value  = "Chelsea"
$ws.Range("D3").Value  = "Huddersfield Town"
$ws.Range("C4").Value  = "Southampton"
$ws.Range("C6").Value  = "Manchester City"
$ws.Range("D6").Value  = "West Ham United"
$ws.Range("C7").Value  = "Wolverhampton"
$ws.Range("C8").Value  = "Leicester City"
$ws.Range("D8").Value  = "Fulham"
$ws.Range("C10").Value = "Crystal Palace"
$ws.Range("C12").Value = "Tottenham Hotspur"
$ws.Range("C14").Value = "Liverpool"
$ws.Range("D13").Value = "Palace"
$ws.Range("F10").Value = "Joker"
$ws.Range("G10").Value = "Based on biggest draw prob and biggest increase in prob"
$ws.Range("C15").Value = "everton"
$ws.Range("C16").Value = "watford"

# ---------------------------------------------------------------------------
# Pass 2: write the real content / formatting for every touched cell.
# ---------------------------------------------------------------------------

# --- Week 1 (row 2) ---
$ws.Range("C2").Value = "Manchester United"
$ws.Range("D2").Value = "Burnley"
$ws.Range("E2").Value = 0.6768

# --- Week 2 (row 3) ---
$ws.Range("C3").Value = "Chelsea"
$ws.Range("D3").Value = "Huddersfield Town"
$ws.Range("E3").Value = 0.785999999999999

# --- Week 3 (row 4) - previously blank styled cells, now plain values ---
$ws.Range("C4:E4").Style = "Normal"
$ws.Range("C4").Value = "Southampton"
$ws.Range("D4").Value = "Cardiff City"
$ws.Range("E4").Value = 0.555

# --- Week 4 (row 6) ---
$ws.Range("C6").Value = "Manchester City"
$ws.Range("D6").Value = "West Ham United"
$ws.Range("E6").Value = 0.8485

# --- Week 5 (row 7) - previously blank styled cells, now plain values ---
$ws.Range("C7:E7").Style = "Normal"
$ws.Range("C7").Value = "Wolverhampton"
$ws.Range("D7").Value = "Cardiff City"
$ws.Range("E7").Value = 0.5626

# --- Week 6 (row 8) - previously blank styled cells, now plain values ---
$ws.Range("C8:E8").Style = "Normal"
$ws.Range("C8").Value = "Leicester City"
$ws.Range("D8").Value = "Fulham"
$ws.Range("E8").Value = 0.5573
# Remove the old "city" note in I8
$ws.Range("I8").ClearContents() | Out-Null

# --- Row 9 (FA Cup Quarters week) - remove the "liverpool" text but keep the styled blank cell ---
$ws.Range("I9").ClearContents() | Out-Null

# --- Week 7 (row 10) - previously blank styled cells, now plain values + Joker note ---
$ws.Range("C10:E10").Style = "Normal"
$ws.Range("C10").Value = "Crystal Palace"
$ws.Range("D10").Value = "Huddersfield Town"
$ws.Range("E10").Value = 0.5783
$ws.Range("F10").Value = "Joker"
$ws.Range("G10").Value = "Based on biggest draw prob and biggest increase in prob"

# --- Week 8 (row 12) ---
$ws.Range("C12").Style = "Normal"
$ws.Range("C12").Value = "Tottenham Hotspur"
$ws.Range("D12").Value = "Huddersfield Town"
$ws.Range("E12").Value = 0.8031

# --- Week 9 (row 13) - C13 keeps its styled (Helvetica) format, D13 gets the same format ---
$ws.Range("C13").Value = "Arsenal"
$ws.Range("C13").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("D13").Value = "Palace"
$ws.Range("E13").Value = 0.5661

# --- Week 10 (row 14) ---
$ws.Range("C14").Value = "Liverpool"
$ws.Range("D14").Value = "Huddersfield Town"
$ws.Range("E14").Value = 0.8921

# --- New row 15: note styled like the other Helvetica cells ---
$ws.Range("C9").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = "everton"

# --- New row 16: plain note ---
$ws.Range("C16").Value = "watford"

# --- New rows 19-28: blank cells formatted like a date column ---
$ws.Range("B2").Copy()
$ws.Range("D19:D28").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Range("C17").Select() | Out-Null
